$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(350, 25, 250, 400, 50, 30, 70)
$startRow = 16

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("J6").Select()
